$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gesamtergebnis")

# The "balance" (Startguthaben/Endsaldo) total cells previously held the
# placeholder text "N/A"; now they should report an actual (zero) total
# like the rest of the totals row.
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0

# Make the Gesamtergebnis sheet the active sheet/selection, matching the
# workbook's last-saved view state.
$ws.Activate()
$ws.Range("E10").Select()
